$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "66.076.30"
Set-TextValue "E2" "  -0.73%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.293.61"
Set-TextValue "E3" "  -0.89%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.00%  "

# Row 5 - BNB
Set-TextValue "D5" "585.20"
Set-TextValue "E5" "  +2.09%  "

# Row 6 - Solana
Set-TextValue "D6" "181.07"
Set-TextValue "E6" "  -0.88%  "

# Row 7 - XRP
Set-TextValue "D7" "0.650"
Set-TextValue "E7" "  +8.12%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.125"
Set-TextValue "E9" "  -3.09%  "

# Row 10 - Toncoin
Set-TextValue "D10" "6.75"
Set-TextValue "E10" "  +1.64%  "

# Row 11 - Cardano
Set-TextValue "E11" "  +0.31%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "3.866.34"

# Row 13 - TRON
Set-TextValue "E13" "  -4.56%  "

# Row 14 - WrappedBTC
Set-TextValue "D14" "66.135.14"
Set-TextValue "E14" "  -0.74%  "

# Row 15 - Avalanche
Set-TextValue "D15" "26.51"
Set-TextValue "E15" "  -2.50%  "

# Row 16 - ShibaInu
Set-TextValue "D16" "0.0000163"
Set-TextValue "E16" "  -2.27%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "3.284.45"
Set-TextValue "E17" "  -1.62%  "

# Row 18 - BitcoinCash
Set-TextValue "D18" "431.42"
Set-TextValue "E18" "  -1.71%  "

# Row 19 - Chainlink
Set-TextValue "D19" "13.26"
Set-TextValue "E19" "  -3.64%  "

# Row 20 - Polkadot
Set-TextValue "E20" "  -3.17%  "

# Row 21 - Uniswap
Set-TextValue "D21" "7.42"
Set-TextValue "E21" "  -3.13%  "

# Row 22 - Litecoin
Set-TextValue "D22" "72.29"
Set-TextValue "E22" "  -2.21%  "

# Row 23 - Dai
Set-TextValue "E23" "  +0.15%  "

# Row 24 - LEO
Set-TextValue "E24" "  +0.33%  "

# Row 25 - WrappedeETH
Set-TextValue "D25" "3.435.08"
Set-TextValue "E25" "  -0.78%  "

# Row 26 - Polygon
Set-TextValue "E26" "  -0.93%  "

# Row 27 - PEPE
Set-TextValue "E27" "  -3.79%  "

# Row 28 - Kaspa
Set-TextValue "E28" "  +1.45%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextValue "E29" "  -1.93%  "

# Row 30 - Binance-PegBSC-USD
Set-TextValue "D30" "1.00"
Set-TextValue "E30" "  +0.33%  "

# Row 31 - PancakeSwap
Set-TextValue "D31" "1.98"
Set-TextValue "E31" "  +0.60%  "

# Row 32 - EthereumClassic
Set-TextValue "D32" "22.37"
Set-TextValue "E32" "  -2.42%  "

# Row 33 - USDe
Set-TextValue "D33" "1.00"
Set-TextValue "E33" "  -0.01%  "

# Row 34 - NEARProtocol
Set-TextValue "E34" "  -3.09%  "

# Row 35 - Aptos
Set-TextValue "E35" "  -2.46%  "

# Row 36 - Fetch.AI
Set-TextValue "D36" "1.19"
Set-TextValue "E36" "  -2.88%  "

# Row 37 - Monero
Set-TextValue "E37" "  -0.89%  "

# Row 38 - ImmutableX
Set-TextValue "E38" "  -5.23%  "

# Row 39 - EnergySwap
Set-TextValue "D39" "26.62"
Set-TextValue "E39" "  -3.29%  "

# Row 40 - Stacks
Set-TextValue "E40" "  -3.76%  "

# Row 41 - Maker
Set-TextValue "D41" "2.788.35"

# Row 42 - Mantle
Set-TextValue "E42" "  -2.43%  "

# Row 43 - Filecoin
Set-TextValue "D43" "4.34"
Set-TextValue "E43" "  -2.96%  "

# Row 44 - OKB
Set-TextValue "D44" "40.10"
Set-TextValue "E44" "  -0.26%  "

# Row 45 - RenderToken
Set-TextValue "D45" "6.01"
Set-TextValue "E45" "  -3.49%  "

# Row 46 - Hedera
Set-TextValue "E46" "  -2.61%  "

# Row 47 - dogwifhat
Set-TextValue "E47" "  -1.25%  "

# Row 48 - was InjectiveProtocol, now Bittensor
Set-TextValue "B48" "Bittensor"
Set-TextValue "C48" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D48" "315.87"
Set-TextValue "E48" "  -1.17%  "

# Row 49 - was Bittensor, now InjectiveProtocol
Set-TextValue "B49" "InjectiveProtocol"
Set-TextValue "C49" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D49" "23.24"
Set-TextValue "E49" "  -4.17%  "

# Row 50 - VeChain
Set-TextValue "E50" "  -2.07%  "
